$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.355.67'
$ws.Range('E2').Value = '  +2.54%  '
$ws.Range('D3').Value = '3.407.25'
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.25%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +9.51%  '
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.64'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.22%  '
$ws.Range('E12').Value = '  +4.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '686.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('E14').Value = '  +4.00%  '
$ws.Range('D15').Value = '3.959.97'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').Value = '69.486.44'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.412.25'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.121'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.23%  '
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.83'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '564.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.27'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = '3.669.12'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.54%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.141'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.51%  '
$ws.Range('D40').Value = '0.0₃0723'
$ws.Range('E40').Value = '  +9.04%  '
$ws.Range('E41').Value = '  +4.13%  '
$ws.Range('E42').Value = '  +3.02%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.339'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('E45').Value = '  +5.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.68'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.93%  '
$ws.Range('E47').Value = '  +1.50%  '
$ws.Range('E48').Value = '  +5.89%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.90'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('E51').Value = '  +2.41%  '
